$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Final target data for rows 16-29 (B,C,D,E,F,G)
# Columns: B=Tipo Doc, C=N Doc, D=Nombre, E=Periodo Mora, F=Valor Mora, G=Salario Basico
$emilsa = "EMILSA ISABEL LOPEZ DE ORDOÑEZ"
$cesar  = "CESAR JOSE ORDOÑEZ LOPEZ"

$rows = @(
  @{R=16; B="CC"; C="26024057"; D=$emilsa; E="2112"; F=36341; G=908526},
  @{R=17; B="CC"; C="73169249"; D=$cesar;  E="2112"; F=36341; G=877803},
  @{R=18; B="CC"; C="26024057"; D=$emilsa; E="2201"; F=36341; G=908526},
  @{R=19; B="CC"; C="73169249"; D=$cesar;  E="2201"; F=36341; G=877803},
  @{R=20; B="CC"; C="26024057"; D=$emilsa; E="2202"; F=36341; G=908526},
  @{R=21; B="CC"; C="73169249"; D=$cesar;  E="2202"; F=36341; G=877803},
  @{R=22; B="CC"; C="26024057"; D=$emilsa; E="2203"; F=36341; G=908526},
  @{R=23; B="CC"; C="73169249"; D=$cesar;  E="2203"; F=36341; G=877803},
  @{R=24; B="CC"; C="26024057"; D=$emilsa; E="2204"; F=36341; G=908526},
  @{R=25; B="CC"; C="73169249"; D=$cesar;  E="2204"; F=36341; G=877803},
  @{R=26; B="CC"; C="26024057"; D=$emilsa; E="2205"; F=36341; G=908526},
  @{R=27; B="CC"; C="73169249"; D=$cesar;  E="2205"; F=36341; G=877803},
  @{R=28; B="CC"; C="26024057"; D=$emilsa; E="2206"; F=26650; G=908526},
  @{R=29; B="CC"; C="73169249"; D=$cesar;  E="2206"; F=25749; G=877803}
)

foreach ($row in $rows) {
  $r = $row.R
  $ws.Cells.Item($r, 2).Value = $row.B
  $ws.Cells.Item($r, 3).Value = $row.C
  $ws.Cells.Item($r, 4).Value = $row.D
  $ws.Cells.Item($r, 5).Value = $row.E
  $ws.Cells.Item($r, 6).Value = $row.F
  $ws.Cells.Item($r, 7).Value = $row.G
}


